# Auto-generated edit script: updates Leve profit-calculator values per the
# upstream "chore: update Sheets via scheduled runner" price-refresh commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 778.6667
$ws.Range("I6").Value = 778.6667
$ws.Range("K6").Value = 2336.0001
$ws.Range("M6").Value = -2224.0001
$ws.Range("H8").Value = 75
$ws.Range("I8").Value = 75
$ws.Range("K8").Value = 225
$ws.Range("M8").Value = -86
$ws.Range("H17").Value = 633286
$ws.Range("J17").Value = 662496.1
$ws.Range("L17").Value = 1987488.3
$ws.Range("N17").Value = -1987824.3
$ws.Range("H19").Value = 1319.8462
$ws.Range("I19").Value = 974.875
$ws.Range("K19").Value = 974.875
$ws.Range("M19").Value = -799.875
$ws.Range("H46").Value = 5101.8887
$ws.Range("I46").Value = 4489.625
$ws.Range("J46").Value = 10000
$ws.Range("K46").Value = 13468.875
$ws.Range("L46").Value = 30000
$ws.Range("M46").Value = -13349.875
$ws.Range("N46").Value = -30238
$ws.Range("H60").Value = 5101.8887
$ws.Range("I60").Value = 4489.625
$ws.Range("J60").Value = 10000
$ws.Range("K60").Value = 13468.875
$ws.Range("L60").Value = 30000
$ws.Range("M60").Value = -12984.875
$ws.Range("N60").Value = -30968
$ws.Range("H76").Value = 3475315.8
$ws.Range("I76").Value = 4447344
$ws.Range("J76").Value = 3785.7144
$ws.Range("K76").Value = 4447344
$ws.Range("L76").Value = 3785.7144
$ws.Range("M76").Value = -4447029
$ws.Range("N76").Value = -4415.7144
$ws.Range("H79").Value = 3475315.8
$ws.Range("I79").Value = 4447344
$ws.Range("J79").Value = 3785.7144
$ws.Range("K79").Value = 4447344
$ws.Range("L79").Value = 3785.7144
$ws.Range("M79").Value = -4446252
$ws.Range("N79").Value = -5969.7144
$ws.Range("H113").Value = 65409.125
$ws.Range("I113").Value = 93885.45
$ws.Range("J113").Value = 2761.2
$ws.Range("K113").Value = 93885.45
$ws.Range("L113").Value = 2761.2
$ws.Range("M113").Value = -90631.45
$ws.Range("N113").Value = -9269.200000000001
$ws.Range("H132").Value = 198226.11
$ws.Range("I132").Value = 213721.73
$ws.Range("J132").Value = 51017.668
$ws.Range("K132").Value = 641165.1900000001
$ws.Range("L132").Value = 153053.004
$ws.Range("M132").Value = -638635.1900000001
$ws.Range("N132").Value = -158113.004
$ws.Range("H135").Value = 1163.4348
$ws.Range("I135").Value = 966.29266
$ws.Range("J135").Value = 2780
$ws.Range("K135").Value = 8696.63394
$ws.Range("L135").Value = 25020
$ws.Range("M135").Value = -6161.63394
$ws.Range("N135").Value = -30090
$ws.Range("H137").Value = 25642342
$ws.Range("I137").Value = 41667696
$ws.Range("J137").Value = 1778.5333
$ws.Range("K137").Value = 125003088
$ws.Range("L137").Value = 5335.5999
$ws.Range("M137").Value = -125000538
$ws.Range("N137").Value = -10435.5999
$ws.Range("H138").Value = 4168940
$ws.Range("I138").Value = 773211.9399999999
$ws.Range("J138").Value = 9262532
$ws.Range("K138").Value = 2319635.82
$ws.Range("L138").Value = 27787596
$ws.Range("M138").Value = -2314495.82
$ws.Range("N138").Value = -27797876
$ws.Range("H141").Value = 1654.9878
$ws.Range("I141").Value = 1035.9166
$ws.Range("J141").Value = 6112.3
$ws.Range("K141").Value = 3107.7498
$ws.Range("L141").Value = 18336.9
$ws.Range("M141").Value = 2072.2502
$ws.Range("N141").Value = -28696.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18343.223
$ws.Range("I32").Value = 2452.6064
$ws.Range("J32").Value = 503007
$ws.Range("K32").Value = 2452.6064
$ws.Range("L32").Value = 503007
$ws.Range("M32").Value = -2165.6064
$ws.Range("N32").Value = -503581
$ws.Range("H61").Value = 1519.7333
$ws.Range("I61").Value = 1044.9412
$ws.Range("J61").Value = 4210.222
$ws.Range("K61").Value = 1044.9412
$ws.Range("L61").Value = 4210.222
$ws.Range("M61").Value = -832.9412
$ws.Range("N61").Value = -4634.222
$ws.Range("H74").Value = 6480.615
$ws.Range("I74").Value = 1625.6666
$ws.Range("J74").Value = 17404.25
$ws.Range("K74").Value = 1625.6666
$ws.Range("L74").Value = 17404.25
$ws.Range("M74").Value = -751.6666
$ws.Range("N74").Value = -19152.25
$ws.Range("H77").Value = 6480.615
$ws.Range("I77").Value = 1625.6666
$ws.Range("J77").Value = 17404.25
$ws.Range("K77").Value = 8128.333000000001
$ws.Range("L77").Value = 87021.25
$ws.Range("M77").Value = -3760.333000000001
$ws.Range("N77").Value = -95757.25
$ws.Range("H122").Value = 2661.6
$ws.Range("I122").Value = 2575.8333
$ws.Range("J122").Value = 3004.6667
$ws.Range("K122").Value = 7727.499899999999
$ws.Range("L122").Value = 9014.000100000001
$ws.Range("M122").Value = -5277.499899999999
$ws.Range("N122").Value = -13914.0001
$ws.Range("H132").Value = 2777.4707
$ws.Range("I132").Value = 2242.1072
$ws.Range("J132").Value = 5275.8335
$ws.Range("K132").Value = 6726.321599999999
$ws.Range("L132").Value = 15827.5005
$ws.Range("M132").Value = -4196.321599999999
$ws.Range("N132").Value = -20887.5005
$ws.Range("H136").Value = 1519.7333
$ws.Range("I136").Value = 1044.9412
$ws.Range("J136").Value = 4210.222
$ws.Range("K136").Value = 3134.8236
$ws.Range("L136").Value = 12630.666
$ws.Range("M136").Value = -584.8235999999997
$ws.Range("N136").Value = -17730.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2696.4285
$ws.Range("I134").Value = 1812.1471
$ws.Range("J134").Value = 4700.8
$ws.Range("K134").Value = 5436.4413
$ws.Range("L134").Value = 14102.4
$ws.Range("M134").Value = -2901.4413
$ws.Range("N134").Value = -19172.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1984.9143
$ws.Range("I31").Value = 1081.5555
$ws.Range("J31").Value = 2941.4119
$ws.Range("K31").Value = 1081.5555
$ws.Range("L31").Value = 2941.4119
$ws.Range("M31").Value = -786.5554999999999
$ws.Range("N31").Value = -3531.4119
$ws.Range("H34").Value = 1984.9143
$ws.Range("I34").Value = 1081.5555
$ws.Range("J34").Value = 2941.4119
$ws.Range("K34").Value = 1081.5555
$ws.Range("L34").Value = 2941.4119
$ws.Range("M34").Value = -879.5554999999999
$ws.Range("N34").Value = -3345.4119
$ws.Range("H58").Value = 1492.9048
$ws.Range("I58").Value = 794.2759
$ws.Range("J58").Value = 3051.3845
$ws.Range("K58").Value = 794.2759
$ws.Range("L58").Value = 3051.3845
$ws.Range("M58").Value = -591.2759
$ws.Range("N58").Value = -3457.3845
$ws.Range("H132").Value = 2820.4524
$ws.Range("I132").Value = 2247.8438
$ws.Range("J132").Value = 4652.8
$ws.Range("K132").Value = 6743.5314
$ws.Range("L132").Value = 13958.4
$ws.Range("M132").Value = -4213.5314
$ws.Range("N132").Value = -19018.4
$ws.Range("H134").Value = 4346.8823
$ws.Range("I134").Value = 3820.6584
$ws.Range("J134").Value = 6504.4
$ws.Range("K134").Value = 11461.9752
$ws.Range("L134").Value = 19513.2
$ws.Range("M134").Value = -8926.975199999999
$ws.Range("N134").Value = -24583.2
$ws.Range("H136").Value = 1492.9048
$ws.Range("I136").Value = 794.2759
$ws.Range("J136").Value = 3051.3845
$ws.Range("K136").Value = 2382.8277
$ws.Range("L136").Value = 9154.1535
$ws.Range("M136").Value = 167.1723000000002
$ws.Range("N136").Value = -14254.1535

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 644.44446
$ws.Range("J17").Value = 687.375
$ws.Range("L17").Value = 2062.125
$ws.Range("N17").Value = -2400.125
$ws.Range("H34").Value = 5465.2
$ws.Range("I34").Value = 1438
$ws.Range("J34").Value = 8150
$ws.Range("K34").Value = 4314
$ws.Range("L34").Value = 24450
$ws.Range("M34").Value = -4230
$ws.Range("N34").Value = -24618
$ws.Range("H39").Value = 9676.471
$ws.Range("J39").Value = 9676.471
$ws.Range("L39").Value = 29029.413
$ws.Range("N39").Value = -29617.413
$ws.Range("H55").Value = 3900
$ws.Range("J55").Value = 4625
$ws.Range("L55").Value = 13875
$ws.Range("N55").Value = -14229
$ws.Range("H107").Value = 437.85715
$ws.Range("I107").Value = 410.83334
$ws.Range("J107").Value = 600
$ws.Range("K107").Value = 1232.50002
$ws.Range("L107").Value = 1800
$ws.Range("M107").Value = 687.4999800000001
$ws.Range("N107").Value = -5640
$ws.Range("H113").Value = 752.0741
$ws.Range("I113").Value = 629.4737
$ws.Range("J113").Value = 818.6286
$ws.Range("K113").Value = 1888.4211
$ws.Range("L113").Value = 2455.8858
$ws.Range("M113").Value = 281.5789
$ws.Range("N113").Value = -6795.8858
$ws.Range("H114").Value = 2411.8
$ws.Range("I114").Value = 1028
$ws.Range("J114").Value = 2757.75
$ws.Range("K114").Value = 3084
$ws.Range("L114").Value = 8273.25
$ws.Range("M114").Value = 170
$ws.Range("N114").Value = -14781.25
$ws.Range("H117").Value = 176.33333
$ws.Range("I117").Value = 176.33333
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 528.99999
$ws.Range("L117").Value = 0
$ws.Range("M117").Value = 2913.00001
$ws.Range("N117").ClearContents()
$ws.Range("H122").Value = 616.3043
$ws.Range("I122").Value = 313.76923
$ws.Range("K122").Value = 2823.92307
$ws.Range("M122").Value = -373.9230699999998
$ws.Range("H131").Value = 1871.258
$ws.Range("I131").Value = 394.2857
$ws.Range("J131").Value = 2302.0417
$ws.Range("K131").Value = 1182.8571
$ws.Range("L131").Value = 6906.125100000001
$ws.Range("M131").Value = 3857.1429
$ws.Range("N131").Value = -16986.1251
$ws.Range("H132").Value = 12821254
$ws.Range("I132").Value = 750
$ws.Range("J132").Value = 23810258
$ws.Range("K132").Value = 6750
$ws.Range("L132").Value = 214292322
$ws.Range("M132").Value = -4220
$ws.Range("N132").Value = -214297382
$ws.Range("H140").Value = 3399.2642
$ws.Range("I140").Value = 3791.4
$ws.Range("J140").Value = 2636.7778
$ws.Range("K140").Value = 11374.2
$ws.Range("L140").Value = 7910.3334
$ws.Range("M140").Value = -6194.200000000001
$ws.Range("N140").Value = -18270.3334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2154.875
$ws.Range("I102").Value = 1989.7368
$ws.Range("J102").Value = 2782.4
$ws.Range("K102").Value = 1989.7368
$ws.Range("L102").Value = 2782.4
$ws.Range("M102").Value = -367.7367999999999
$ws.Range("N102").Value = -6026.4
$ws.Range("H113").Value = 2036.9048
$ws.Range("I113").Value = 1729.1111
$ws.Range("K113").Value = 1729.1111
$ws.Range("M113").Value = 440.8888999999999
$ws.Range("H132").Value = 2034.1342
$ws.Range("I132").Value = 1830.4844
$ws.Range("K132").Value = 5491.4532
$ws.Range("M132").Value = -2961.4532

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3379.476
$ws.Range("I7").Value = 3068.4285
$ws.Range("J7").Value = 3535
$ws.Range("K7").Value = 3068.4285
$ws.Range("L7").Value = 3535
$ws.Range("M7").Value = -2956.4285
$ws.Range("N7").Value = -3759
$ws.Range("H40").Value = 4493.75
$ws.Range("I40").Value = 3080
$ws.Range("J40").Value = 5136.364
$ws.Range("K40").Value = 3080
$ws.Range("L40").Value = 5136.364
$ws.Range("M40").Value = -2944
$ws.Range("N40").Value = -5408.364
$ws.Range("H42").Value = 17618.6
$ws.Range("J42").Value = 9523.25
$ws.Range("L42").Value = 9523.25
$ws.Range("N42").Value = -10649.25
$ws.Range("H49").Value = 17618.6
$ws.Range("J49").Value = 9523.25
$ws.Range("L49").Value = 9523.25
$ws.Range("N49").Value = -9817.25
$ws.Range("H61").Value = 8178.8486
$ws.Range("J61").Value = 6199.8887
$ws.Range("L61").Value = 6199.8887
$ws.Range("N61").Value = -6603.8887
$ws.Range("H68").Value = 2565.6667
$ws.Range("I68").Value = 2334.3333
$ws.Range("K68").Value = 2334.3333
$ws.Range("M68").Value = -1585.3333
$ws.Range("H71").Value = 2565.6667
$ws.Range("I71").Value = 2334.3333
$ws.Range("K71").Value = 11671.6665
$ws.Range("M71").Value = -7927.666499999999
$ws.Range("H82").Value = 1056.2354
$ws.Range("I82").Value = 933.1667
$ws.Range("J82").Value = 1123.3636
$ws.Range("K82").Value = 933.1667
$ws.Range("L82").Value = 1123.3636
$ws.Range("M82").Value = -572.1667
$ws.Range("N82").Value = -1845.3636
$ws.Range("H85").Value = 1056.2354
$ws.Range("I85").Value = 933.1667
$ws.Range("J85").Value = 1123.3636
$ws.Range("K85").Value = 933.1667
$ws.Range("L85").Value = 1123.3636
$ws.Range("M85").Value = 314.8333
$ws.Range("N85").Value = -3619.3636
$ws.Range("H113").Value = 8178.8486
$ws.Range("J113").Value = 6199.8887
$ws.Range("L113").Value = 6199.8887
$ws.Range("N113").Value = -10539.8887
$ws.Range("H122").Value = 3420
$ws.Range("I122").Value = 2800
$ws.Range("J122").Value = 3833.3333
$ws.Range("K122").Value = 8400
$ws.Range("L122").Value = 11499.9999
$ws.Range("M122").Value = -5950
$ws.Range("N122").Value = -16399.9999
$ws.Range("H126").Value = 3379.476
$ws.Range("I126").Value = 3068.4285
$ws.Range("J126").Value = 3535
$ws.Range("K126").Value = 9205.2855
$ws.Range("L126").Value = 10605
$ws.Range("M126").Value = -6735.2855
$ws.Range("N126").Value = -15545
$ws.Range("H132").Value = 2554.0393
$ws.Range("I132").Value = 1871.9706
$ws.Range("K132").Value = 5615.9118
$ws.Range("M132").Value = -3085.9118
$ws.Range("H136").Value = 2978.9678
$ws.Range("I136").Value = 1769.7
$ws.Range("K136").Value = 5309.1
$ws.Range("M136").Value = -2759.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 39633.332
$ws.Range("J63").Value = 39633.332
$ws.Range("L63").Value = 39633.332
$ws.Range("N63").Value = -40881.332
$ws.Range("H64").Value = 29057
$ws.Range("J64").Value = 29057
$ws.Range("L64").Value = 29057
$ws.Range("N64").Value = -29553
$ws.Range("H66").Value = 39633.332
$ws.Range("J66").Value = 39633.332
$ws.Range("L66").Value = 118899.996
$ws.Range("N66").Value = -125139.996
$ws.Range("H67").Value = 29057
$ws.Range("J67").Value = 29057
$ws.Range("L67").Value = 29057
$ws.Range("N67").Value = -30773
$ws.Range("H113").Value = 564.86664
$ws.Range("I113").Value = 386
$ws.Range("K113").Value = 1158
$ws.Range("M113").Value = 1012
$ws.Range("H132").Value = 9261240
$ws.Range("I132").Value = 12501640
$ws.Range("K132").Value = 37504920
$ws.Range("M132").Value = -37502390
$ws.Range("H136").Value = 7961619
$ws.Range("I136").Value = 10449029
$ws.Range("J136").Value = 1906.5
$ws.Range("K136").Value = 31347087
$ws.Range("L136").Value = 5719.5
$ws.Range("M136").Value = -31344537
$ws.Range("N136").Value = -10819.5
